$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 2.7
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 7.6
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.98
$ws.Range("S2").Value = 1.37
$ws.Range("T2").Value = 2.85
$ws.Range("U2").Value = 1.6
$ws.Range("V2").Value = 2.22
$ws.Range("X2").Value = 12
$ws.Range("AB2").Value = 22
$ws.Range("AC2").Value = 7.6
$ws.Range("AD2").Value = 6.3
$ws.Range("AE2").Value = 11.75
$ws.Range("AI2").Value = 18
$ws.Range("AL2").Value = 27
$ws.Range("AO2").Value = 11
$ws.Range("AP2").Value = 16.5
$ws.Range("AR2").Value = 60
$ws.Range("AT2").Value = 2.85
$ws.Range("AW2").Value = 5.3
$ws.Range("AX2").Value = 18.5

$wb.Save()
